$d = $word.ActiveDocument

$replacements = @(
    @{old = "349×6="; new = "725×6="},
    @{old = "410×6="; new = "117×9="},
    @{old = "986×9="; new = "755×7="},
    @{old = "919×7="; new = "742×2="},
    @{old = "868×4="; new = "892×6="},
    @{old = "346×6="; new = "835×8="},
    @{old = "947×5="; new = "770×3="},
    @{old = "945×4="; new = "749×7="},
    @{old = "412×2="; new = "300×2="},
    @{old = "547×7="; new = "929×9="},
    @{old = "869×5="; new = "855×9="},
    @{old = "820×6="; new = "152×4="},
    @{old = "315×5="; new = "700×7="},
    @{old = "335×8="; new = "350×9="},
    @{old = "501×6="; new = "312×7="},
    @{old = "112×3="; new = "381×9="},
    @{old = "695×2="; new = "300×5="},
    @{old = "263×6="; new = "718×8="},
    @{old = "814×7="; new = "233×5="},
    @{old = "410×5="; new = "102×9="},
    @{old = "917×6="; new = "334×7="},
    @{old = "444×3="; new = "683×7="},
    @{old = "568×4="; new = "658×7="},
    @{old = "396×7="; new = "247×6="},
    @{old = "980×4="; new = "659×6="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
